$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 5 = "BIOS/POST/SETUP/CMOS/BATERIA/BOOT:MBR;GPT"
# Paragraph 6 = "<tab>https://homepages.dcc.ufmg.br/~cesarfmc/classes/manut2/TeoricaSetup.pdf "
# Target: the BIOS line gets reworded/re-colored and the two paragraphs become one
# (the tab+link paragraph is absorbed into the BIOS paragraph).

$url = "https://homepages.dcc.ufmg.br/~cesarfmc/classes/manut2/TeoricaSetup.pdf"
$newText = "BIOS/POST/SETUP/CMOS/ bootloader/BOOT:MBR;GPT; " + $url + " "

$para5 = $tr.Paragraphs(5, 1)
$base = $para5.Start
$para5.Text = $newText

function Set-Run($start, $length, $rgb) {
    $rng = $tr.Characters($start, $length)
    $rng.Font.Bold = -1
    if ($rgb -ne $null) {
        $rng.Font.Color.RGB = $rgb
    }
}

# "BIOS" - bold, color 002060 (COM RGB is stored BGR, so 0x602000 round-trips to 002060)
Set-Run ($base + 0) 4 0x602000
# "/POST/SETUP/CMOS/" - bold, default text color
Set-Run ($base + 4) 17 $null
# " " - bold, red FF0000 (BGR 0x0000FF)
Set-Run ($base + 21) 1 0x0000FF
# "bootloader" - bold, default text color
Set-Run ($base + 22) 10 $null
# "/BOOT" - bold, red FF0000
Set-Run ($base + 32) 5 0x0000FF
# ":" - bold, default text color
Set-Run ($base + 37) 1 $null
# "MBR;GPT; " - bold, red FF0000
Set-Run ($base + 38) 9 0x0000FF

# The URL keeps the non-bold plain style from the original link line and
# gets its hyperlink re-attached (re-using the rId3 relationship, since the
# target address already matches it).
$urlStart = $base + 47
$urlRange = $tr.Characters($urlStart, $url.Length)
$urlRange.Font.Bold = 0
$hl = $urlRange.ActionSettings(1).Hyperlink
$hl.Address = $url

# Trailing space after the URL, plain style (matches original run).
$trailRange = $tr.Characters($urlStart + $url.Length, 1)
$trailRange.Font.Bold = 0

# Absorb the old tab+link paragraph into the rebuilt paragraph above.
$para6 = $tr.Paragraphs(6, 1)
$para6.Delete()
